# Refresh the cryptos table: updated prices (column D) and 1h volume
# percentage changes (column E) for every coin row, plus two row pairs
# (38/39 and 40/41) whose coin name, link, price and change all swapped
# places because the ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D prices are stored as literal text (e.g. '1.000', '27.497.82')
# and would otherwise be auto-converted to numbers by Excel (losing the
# trailing zeros / thousands-dot formatting), so we prefix them with a
# leading apostrophe - exactly like typing '1.000 into the cell - to force
# literal text entry.

$ws.Range('D2').Value = "'" + '27.497.82'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').Value = "'" + '1.749.76'
$ws.Range('E3').Value = '  -2.39%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'" + '324.14'
$ws.Range('D6').Value = "'" + '1.000'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = "'" + '0.4457'
$ws.Range('E7').Value = '  +4.05%  '
$ws.Range('D8').Value = "'" + '0.3605'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = "'" + '0.07495'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = "'" + '42.01'
$ws.Range('E10').Value = '  -6.03%  '
$ws.Range('D11').Value = "'" + '1.092'
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('D12').Value = "'" + '0.9998'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = "'" + '20.61'
$ws.Range('E13').Value = '  -4.67%  '
$ws.Range('D14').Value = "'" + '6.020'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = "'" + '7.129'
$ws.Range('E15').Value = '  -2.62%  '
$ws.Range('D16').Value = "'" + '1.754.77'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').Value = "'" + '92.58'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = "'" + '0.00001060'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').Value = "'" + '0.06404'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = "'" + '1.000'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = "'" + '16.82'
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').Value = "'" + '5.849'
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').Value = "'" + '27.547.89'
$ws.Range('E23').Value = '  -2.16%  '
$ws.Range('D24').Value = "'" + '11.15'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = "'" + '2.096'
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('D26').Value = "'" + '161.89'
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('D27').Value = "'" + '20.47'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').Value = "'" + '1.950.68'
$ws.Range('E28').Value = '  -3.09%  '
$ws.Range('D29').Value = "'" + '2.083'
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('D30').Value = "'" + '124.80'
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('D31').Value = "'" + '1.084'
$ws.Range('E31').Value = '  -7.24%  '
$ws.Range('E32').Value = '  +3.77%  '
$ws.Range('D33').Value = "'" + '0.09009'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').Value = "'" + '5.518'
$ws.Range('E34').Value = '  -4.84%  '
$ws.Range('D35').Value = "'" + '12.01'
$ws.Range('E35').Value = '  -5.37%  '
$ws.Range('D36').Value = "'" + '0.02299'
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('D37').Value = "'" + '0.2084'
$ws.Range('E37').Value = '  -1.38%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'" + '0.06003'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = "'" + '0.6339'
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = "'" + '4.940'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'" + '1.205'
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('D42').Value = "'" + '0.9991'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = "'" + '1.385'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').Value = "'" + '7.779'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').Value = "'" + '13.12'
$ws.Range('E45').Value = '  -3.15%  '
$ws.Range('D46').Value = "'" + '3.710'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = "'" + '0.5882'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').Value = "'" + '121.32'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('D49').Value = "'" + '1.950'
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('D50').Value = "'" + '1.150'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').Value = "'" + '0.06858'
$ws.Range('E51').Value = '  -1.65%  '
